$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C4").Value = -12.481
$ws.Range("C7").Value = -12.793
$ws.Range("D7").Value = -7.7
$ws.Range("D15").Value = -8.24
$ws.Range("C16").Value = -13.331
$ws.Range("E16").Value = 16.65
$ws.Range("E19").Value = 16.45
$ws.Range("D21").Value = -8.1
$ws.Range("D22").Value = -7.900000000000001
$ws.Range("D23").Value = -7.856
$ws.Range("C28").Value = -13.219
$ws.Range("C29").Value = -11.997
$ws.Range("C32").Value = -13.434
$ws.Range("D34").Value = -7.782999999999999
$ws.Range("E36").Value = 16.744
$ws.Range("C40").Value = -12.644
$ws.Range("D43").Value = -7.707000000000001
$ws.Range("D45").Value = -7.703999999999999
$ws.Range("E46").Value = 16.773
$ws.Range("D50").Value = -8.217000000000002
$ws.Range("E50").Value = 16.647
$ws.Range("D51").Value = -8.349
$ws.Range("C52").Value = -11.614
$ws.Range("C57").Value = -13.622
$ws.Range("C66").Value = -11.527
$ws.Range("D66").Value = -7.447
$ws.Range("D67").Value = -6.768000000000001
$ws.Range("D79").Value = -7.703
$ws.Range("D84").Value = -8.141000000000002
$ws.Range("D92").Value = -6.552
$ws.Range("E95").Value = 17.52
$ws.Range("D97").Value = -8.184000000000001
$ws.Range("E97").Value = 16.597
$ws.Range("C100").Value = -13.172
